# Append 9 new NBA game rows (888-896) to Sheet1 and correct the overtime
# flag on the last existing row (887), matching the upstream data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 887 already exists (Washington Wizards @ Los Angeles Lakers) but was
# recorded as an overtime game ("Yes"); the refreshed source marks it "OT".
$ws.Cells.Item(887, 5).Value = "OT"

# New rows appended at the bottom of the table.
$newRows = @(
    ,@("Cleveland Cavaliers",     110, "Detroit Pistons",         100, "No", 17832, "Little Caesars Arena", "Cleveland Cavaliers",    "Detroit Pistons")
    ,@("Charlotte Hornets",       114, "Philadelphia 76ers",      121, "No", 17832, "Wells Fargo Center",   "Philadelphia 76ers",     "Charlotte Hornets")
    ,@("Dallas Mavericks",        110, "Boston Celtics",          138, "No", 17832, "TD Garden",            "Boston Celtics",         "Dallas Mavericks")
    ,@("Golden State Warriors",   120, "Toronto Raptors",         105, "No", 17832, "Scotiabank Arena",     "Golden State Warriors",  "Toronto Raptors")
    ,@("Portland Trail Blazers",  122, "Memphis Grizzlies",        92, "No", 17832, "FedEx Forum",          "Portland Trail Blazers", "Memphis Grizzlies")
    ,@("Sacramento Kings",        124, "Minnesota Timberwolves",  120, "OT", 17832, "Target Center",        "Sacramento Kings",       "Minnesota Timberwolves")
    ,@("Indiana Pacers",          102, "New Orleans Pelicans",    129, "No", 17832, "Smoothie King Center", "New Orleans Pelicans",   "Indiana Pacers")
    ,@("Milwaukee Bucks",         113, "Chicago Bulls",            97, "No", 17832, "United Center",        "Milwaukee Bucks",        "Chicago Bulls")
    ,@("Washington Wizards",      115, "Los Angeles Clippers",    140, "No", 17832, "Crypto.com Arena",     "Los Angeles Clippers",   "Washington Wizards")
)

$startRow = 888
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $values = $newRows[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $values[$c]
    }
}

# Match the author's final view state: scrolled down with D896 selected.
[void]$ws.Range("D896").Select()
$excel.ActiveWindow.ScrollRow = 865
